$wb = $excel.ActiveWorkbook

# 1. Reorder sheets: put "review_info" before "hotel_info"
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))

# 2. Insert a new "State" column into "hotel_info" right after "Hotel_Name" (i.e. before "City")
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Columns.Item(3).Insert()

$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"
